$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38, shifting existing rows 38-88 down to 39-89
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new record's data
$ws.Cells.Item(38, 1).Value = 1
$ws.Cells.Item(38, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(38, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(38, 4).Value = 44757
$ws.Cells.Item(38, 5).Value = 15
$ws.Cells.Item(38, 6).Value = 100112021
$ws.Cells.Item(38, 7).Value = "Ají"
$ws.Cells.Item(38, 8).Value = "Inferno"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 140
$ws.Cells.Item(38, 11).Value = 10000
$ws.Cells.Item(38, 12).Value = 11000
$ws.Cells.Item(38, 13).Value = 10500
$ws.Cells.Item(38, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(38, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(38, 16).Value = 700
$ws.Cells.Item(38, 17).Value = 15
$ws.Cells.Item(38, 18).Value = "Hortaliza"
